# Apply updated crypto price/volume figures (and restore original B/C/D/E
# ordering swap for rows 29-30) as captured in the authoritative XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.627.59"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3
$ws.Range("D3").Value = "1.844.24"
$ws.Range("E3").Value = "  -0.52%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'315.88"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.4230"
$ws.Range("E7").Value = "  -2.94%  "
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "  -1.00%  "

# Row 9
$ws.Range("D9").Value = "'45.20"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.07265"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.8904"
$ws.Range("E11").Value = "  -4.78%  "
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'20.64"
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "1.861.62"
$ws.Range("E13").Value = "  +0.36%  "

# Row 14
$ws.Range("D14").Value = "'6.572"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'5.336"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.06872"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "  -0.43%  "

# Row 18
$ws.Range("D18").Value = "'78.83"
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.000008842"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "  -2.54%  "

# Row 22
$ws.Range("D22").Value = "27.606.98"
$ws.Range("E22").Value = "  -1.05%  "

# Row 23
$ws.Range("D23").Value = "'4.978"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'10.55"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "2.090.75"
$ws.Range("E25").Value = "  -0.40%  "

# Row 26
$ws.Range("D26").Value = "'1.945"
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'154.72"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "  +0.89%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'119.55"
$ws.Range("E29").Value = "  +5.91%  "
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'5.234"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "  +6.61%  "

# Row 32
$ws.Range("D32").Value = "'0.08915"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'0.7788"
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'4.561"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'2.944"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "  -6.12%  "

# Row 37
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.05395"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'1.096"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.01925"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D40").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'6.849"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.5059"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = "  -1.78%  "

# Row 45
$ws.Range("D45").Value = "'8.268"
$ws.Range("E45").Value = "  -5.05%  "
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "  -1.34%  "

# Row 47
$ws.Range("D47").Value = "'10.36"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.4698"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'104.42"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "  -2.41%  "
